# Update BunkerPrices at 2025-03-19 14:59
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Montevideo/New York header labels (columns Y and Z) on row 1
$ws.Range("Y1").Value = "New York"
$ws.Range("Z1").Value = "Montevideo"

# Swap the corresponding data values on row 2
$ws.Range("Y2").Value = 534
$ws.Range("Z2").Value = 552

# Swap the corresponding data values on row 3
$ws.Range("Y3").Value = 535
$ws.Range("Z3").Value = 553

# AH3's date now uses the same number format as AH2 (YYYY-MM-DD HH:MM:SS)
$ws.Range("AH3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append a brand new data row (row 4) with a full set of bunker prices
$ws.Range("A4").Value = 577
$ws.Range("B4").Value = 646
$ws.Range("C4").Value = 499
$ws.Range("D4").Value = 571
$ws.Range("E4").Value = 620
$ws.Range("F4").Value = 645
$ws.Range("G4").Value = 502
$ws.Range("H4").Value = 516
$ws.Range("I4").Value = 563
$ws.Range("J4").Value = 520
$ws.Range("K4").Value = 582
$ws.Range("L4").Value = 516
$ws.Range("M4").Value = 532
$ws.Range("N4").Value = 883
$ws.Range("O4").Value = 578
$ws.Range("P4").Value = 526
$ws.Range("Q4").Value = 502
$ws.Range("R4").Value = 530
$ws.Range("S4").Value = 605
$ws.Range("T4").Value = 646
$ws.Range("U4").Value = 588
$ws.Range("V4").Value = 490
$ws.Range("W4").Value = 560
$ws.Range("X4").Value = 528
$ws.Range("Y4").Value = 530
$ws.Range("Z4").Value = 552
$ws.Range("AA4").Value = 502
$ws.Range("AB4").Value = 543
$ws.Range("AC4").Value = 584.5
$ws.Range("AD4").Value = 515
$ws.Range("AE4").Value = 514
$ws.Range("AF4").Value = 527
$ws.Range("AG4").Value = 490
$ws.Range("AH4").Value = 45730
$ws.Range("AH4").NumberFormat = "YYYY-MM-DD"
$ws.Range("AI4").Value = 505
$ws.Range("AJ4").Value = 537
$ws.Range("AK4").Value = 520
$ws.Range("AL4").Value = 750
$ws.Range("AM4").Value = 646
$ws.Range("AN4").Value = 610
$ws.Range("AO4").Value = 500
$ws.Range("AP4").Value = 640
$ws.Range("AQ4").Value = 760
$ws.Range("AR4").Value = 515
$ws.Range("AS4").Value = 499
$ws.Range("AT4").Value = 566
$ws.Range("AU4").Value = 578
$ws.Range("AV4").Value = 637
